# Scheduled runner: refresh Universalis market-price snapshot columns
# (currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ)
# across the Anima_Profits crafting-job sheets (cols H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 581.9286
$ws.Range("I33").Value = 664.6667
$ws.Range("J33").Value = 85.5
$ws.Range("K33").Value = 664.6667
$ws.Range("L33").Value = 85.5
$ws.Range("M33").Value = -435.6667
$ws.Range("N33").Value = -543.5

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 431.9091
$ws.Range("I125").Value = 533.3333
$ws.Range("J125").Value = 393.875
$ws.Range("K125").Value = 4799.9997
$ws.Range("L125").Value = 3544.875
$ws.Range("M125").Value = -2339.9997
$ws.Range("N125").Value = -8464.875

# Row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 1923.5
$ws.Range("I131").Value = 443
$ws.Range("J131").Value = 2596.4546
$ws.Range("K131").Value = 1329
$ws.Range("L131").Value = 7789.3638
$ws.Range("M131").Value = 3711
$ws.Range("N131").Value = -17869.3638

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1335.1333
$ws.Range("I137").Value = 969.2632
$ws.Range("J137").Value = 1967.091
$ws.Range("K137").Value = 2907.7896
$ws.Range("L137").Value = 5901.272999999999
$ws.Range("M137").Value = -357.7896000000001
$ws.Range("N137").Value = -11001.273

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3499.8262
$ws.Range("I138").Value = 3637.5293
$ws.Range("J138").Value = 3454.8076
$ws.Range("K138").Value = 10912.5879
$ws.Range("L138").Value = 10364.4228
$ws.Range("M138").Value = -5772.5879
$ws.Range("N138").Value = -20644.4228

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 3884.4092
$ws.Range("I141").Value = 1974.8823
$ws.Range("K141").Value = 5924.6469
$ws.Range("M141").Value = -744.6468999999997

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2722.8865
$ws.Range("I74").Value = 1608.3478
$ws.Range("J74").Value = 3943.5715
$ws.Range("K74").Value = 1608.3478
$ws.Range("L74").Value = 3943.5715
$ws.Range("M74").Value = -734.3478
$ws.Range("N74").Value = -5691.5715

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2722.8865
$ws.Range("I77").Value = 1608.3478
$ws.Range("J77").Value = 3943.5715
$ws.Range("K77").Value = 8041.739
$ws.Range("L77").Value = 19717.8575
$ws.Range("M77").Value = -3673.739
$ws.Range("N77").Value = -28453.8575

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 4774.3335
$ws.Range("I132").Value = 3551.7693
$ws.Range("K132").Value = 10655.3079
$ws.Range("M132").Value = -8125.3079

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2489.7742
$ws.Range("I86").Value = 2217.2
$ws.Range("J86").Value = 2542.1924
$ws.Range("K86").Value = 2217.2
$ws.Range("L86").Value = 2542.1924
$ws.Range("M86").Value = -1094.2
$ws.Range("N86").Value = -4788.1924

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2489.7742
$ws.Range("I89").Value = 2217.2
$ws.Range("J89").Value = 2542.1924
$ws.Range("K89").Value = 11086
$ws.Range("L89").Value = 12710.962
$ws.Range("M89").Value = -5470
$ws.Range("N89").Value = -23942.962

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2149.4666
$ws.Range("I134").Value = 1808.2368
$ws.Range("K134").Value = 5424.7104
$ws.Range("M134").Value = -2889.7104

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 6222.709
$ws.Range("I31").Value = 1338.1428
$ws.Range("J31").Value = 9239.647000000001
$ws.Range("K31").Value = 1338.1428
$ws.Range("L31").Value = 9239.647000000001
$ws.Range("M31").Value = -1043.1428
$ws.Range("N31").Value = -9829.647000000001

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 6222.709
$ws.Range("I34").Value = 1338.1428
$ws.Range("J34").Value = 9239.647000000001
$ws.Range("K34").Value = 1338.1428
$ws.Range("L34").Value = 9239.647000000001
$ws.Range("M34").Value = -1136.1428
$ws.Range("N34").Value = -9643.647000000001

# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 1974.4166
$ws.Range("I99").Value = 1571.4286
$ws.Range("J99").Value = 2071.6897
$ws.Range("K99").Value = 1571.4286
$ws.Range("L99").Value = 2071.6897
$ws.Range("M99").Value = -73.42859999999996
$ws.Range("N99").Value = -5067.6897

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 1683.1666
$ws.Range("I122").Value = 1366.3334
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4099.0002
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1649.0002
$ws.Range("N122").Value = -10900

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 1974.4166
$ws.Range("I126").Value = 1571.4286
$ws.Range("J126").Value = 2071.6897
$ws.Range("K126").Value = 4714.2858
$ws.Range("L126").Value = 6215.0691
$ws.Range("M126").Value = -2244.2858
$ws.Range("N126").Value = -11155.0691

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 5210328.5
$ws.Range("I132").Value = 1777.3
$ws.Range("K132").Value = 5331.9
$ws.Range("M132").Value = -2801.9

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 6433.926
$ws.Range("I134").Value = 6441.4546
$ws.Range("J134").Value = 6400.8
$ws.Range("K134").Value = 19324.3638
$ws.Range("L134").Value = 19202.4
$ws.Range("M134").Value = -16789.3638
$ws.Range("N134").Value = -24272.4

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 646.1622
$ws.Range("I5").Value = 467.2
$ws.Range("J5").Value = 1413.1428
$ws.Range("K5").Value = 1401.6
$ws.Range("L5").Value = 4239.428400000001
$ws.Range("M5").Value = -1289.6
$ws.Range("N5").Value = -4463.428400000001

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 2746.4583
$ws.Range("I122").Value = 418.41666
$ws.Range("K122").Value = 3765.74994
$ws.Range("M122").Value = -1315.74994

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 646.1622
$ws.Range("I135").Value = 467.2
$ws.Range("J135").Value = 1413.1428
$ws.Range("K135").Value = 4204.8
$ws.Range("L135").Value = 12718.2852
$ws.Range("M135").Value = -1669.8
$ws.Range("N135").Value = -17788.2852

$ws = $wb.Worksheets.Item("GSM")
# Row 43: Get the Green Stuff / Malachite Earrings
$ws.Range("H43").Value = 12918.8
$ws.Range("I43").Value = 2298.5
$ws.Range("K43").Value = 2298.5
$ws.Range("M43").Value = -2147.5

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 2351002.8
$ws.Range("I80").Value = 4500502.5
$ws.Range("J80").Value = 201503
$ws.Range("K80").Value = 4500502.5
$ws.Range("L80").Value = 201503
$ws.Range("M80").Value = -4499504.5
$ws.Range("N80").Value = -203499

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 2351002.8
$ws.Range("I83").Value = 4500502.5
$ws.Range("J83").Value = 201503
$ws.Range("K83").Value = 22502512.5
$ws.Range("L83").Value = 1007515
$ws.Range("M83").Value = -22497520.5
$ws.Range("N83").Value = -1017499

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3466.2
$ws.Range("I122").Value = 3133.6667
$ws.Range("J122").Value = 3965
$ws.Range("K122").Value = 9401.000100000001
$ws.Range("L122").Value = 11895
$ws.Range("M122").Value = -6951.000100000001
$ws.Range("N122").Value = -16795

# Row 140: Worqor Zormor or Bust / Gargantuaskin Shoes of Healing
$ws.Range("H140").Value = 67369.60000000001
$ws.Range("J140").Value = 67369.60000000001
$ws.Range("L140").Value = 67369.60000000001
$ws.Range("N140").Value = -77729.60000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 54: No Country for Cold Men / Woolen Tights
$ws.Range("H54").Value = 8432
$ws.Range("J54").Value = 8432
$ws.Range("L54").Value = 8432
$ws.Range("N54").Value = -9472

# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 4457.8667
$ws.Range("I81").Value = 5258.75
$ws.Range("J81").Value = 3542.5715
$ws.Range("K81").Value = 10517.5
$ws.Range("L81").Value = 7085.143
$ws.Range("M81").Value = -9456.5
$ws.Range("N81").Value = -9207.143

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 4457.8667
$ws.Range("I84").Value = 5258.75
$ws.Range("J84").Value = 3542.5715
$ws.Range("K84").Value = 52587.5
$ws.Range("L84").Value = 35425.715
$ws.Range("M84").Value = -47283.5
$ws.Range("N84").Value = -46033.715

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1491
$ws.Range("I113").Value = 1649.5883
$ws.Range("J113").Value = 1105.8572
$ws.Range("K113").Value = 4948.7649
$ws.Range("L113").Value = 3317.5716
$ws.Range("M113").Value = -2778.7649
$ws.Range("N113").Value = -7657.571599999999
